$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values: force text to preserve original string formatting
# (these values look numeric and Excel would otherwise coerce them to floats,
# changing the stored representation away from the literal text in the source).
foreach ($pair in @(
        @("D2", "28.097.37"),
        @("D3", "1.863.60"),
        @("D4", "1.006"),
        @("D5", "312.10"),
        @("D7", "0.4988"),
        @("D8", "0.3909"),
        @("D9", "0.09619"),
        @("D10", "1.129"),
        @("D11", "40.88"),
        @("D12", "6.441"),
        @("D13", "20.80"),
        @("D14", "1.871.42"),
        @("D15", "1.006"),
        @("D16", "7.365"),
        @("D17", "0.00001124"),
        @("D18", "92.94"),
        @("D19", "0.06608"),
        @("D20", "1.003"),
        @("D21", "17.34"),
        @("D22", "6.124"),
        @("D23", "28.195.02"),
        @("D24", "11.23"),
        @("D25", "2.282"),
        @("D26", "2.538"),
        @("D27", "2.088.36"),
        @("D28", "21.10"),
        @("D29", "157.70"),
        @("D30", "127.41"),
        @("D31", "0.1058"),
        @("D32", "1.051"),
        @("D33", "5.598"),
        @("D34", "3.613"),
        @("D35", "0.06734"),
        @("D36", "9.457"),
        @("D37", "0.02382"),
        @("D38", "0.2171"),
        @("D39", "4.994"),
        @("D40", "11.44"),
        @("D41", "0.6249"),
        @("D42", "1.172"),
        @("D44", "13.54"),
        @("D45", "0.5968"),
        @("D46", "3.668"),
        @("D48", "123.94"),
        @("D49", "1.978"),
        @("D50", "1.192"),
        @("D51", "0.06834")
    )) {
    $addr = $pair[0]
    $val = $pair[1]
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Other columns (B, C, E): plain text values, no numeric coercion risk
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  -2.37%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  +24.38%  "
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("E14").Value = "  +4.26%  "
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  +4.86%  "
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +4.95%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E28").Value = "  +4.27%  "
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  -2.81%  "
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("E36").Value = "  +4.72%  "
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  +1.52%  "
